$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 2.619789
$ws.Range("D6").Value = 1.993425
$ws.Range("D7").Value = 2.797465
$ws.Range("D8").Value = 2.468092
$ws.Range("D9").Value = 1.862469
$ws.Range("D10").Value = 1.488834
$ws.Range("D11").Value = 1.048059
$ws.Range("D12").Value = 1.135074
$ws.Range("D13").Value = 1.947028
$ws.Range("D14").Value = 1.327318
$ws.Range("D15").Value = 1.809966
$ws.Range("D16").Value = 2.64012
$ws.Range("D24").Value = 2.051499
$ws.Range("D25").Value = 1.513867
$ws.Range("D26").Value = 1.091827
$ws.Range("D30").Value = 4.607615
$ws.Range("D31").Value = 4.748995
$ws.Range("D36").Value = 6.238209
$ws.Range("D37").Value = 6.188049
$ws.Range("D38").Value = 6.470884
$ws.Range("D39").Value = 6.755089
$ws.Range("D42").Value = 6.294948
$ws.Range("D43").Value = 6.968059
$ws.Range("D44").Value = 7.427834
$ws.Range("D45").Value = 7.174165
$ws.Range("D46").Value = 5.541945
$ws.Range("D47").Value = 5.629699
$ws.Range("D48").Value = 4.295833
$ws.Range("D49").Value = 5.436639
$ws.Range("D50").Value = 5.840309
$ws.Range("D51").Value = 6.279081
$ws.Range("D52").Value = 4.190528
$ws.Range("D53").Value = 4.804808
$ws.Range("D54").Value = 4.085222
$ws.Range("D55").Value = 4.050121
$ws.Range("D56").Value = 4.295833
$ws.Range("D57").Value = 3.444616
$ws.Range("D58").Value = 3.549921
$ws.Range("D59").Value = 3.843898
$ws.Range("D60").Value = 4.042422

$ws.Range("D7").Select()
